# Insert two new rows at 461-462 (everything from the old row 461 onward
# shifts down by two rows, old row 461 -> new row 463, ..., old row 555 ->
# new row 557).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("461:462").Insert()

# New row 461
$ws.Range("A461").Value = 11
$ws.Range("B461").Value = "Vega Monumental Concepción"
$ws.Range("C461").Value = "Bíobío"
$ws.Range("D461").Value = 45211
$ws.Range("E461").Value = 8
$ws.Range("F461").Value = "Fruta"
$ws.Range("G461").Value = 100102
$ws.Range("H461").Value = "Cítricos"
$ws.Range("I461").Value = 100102005
$ws.Range("J461").Value = "Naranja"
$ws.Range("K461").Value = "Lane Late"
$ws.Range("L461").Value = "Primera"
$ws.Range("M461").Value = 200
$ws.Range("N461").Value = 9000
$ws.Range("O461").Value = 10000
$ws.Range("P461").Value = 9500
$ws.Range("Q461").Value = "$/bandeja 15 kilos granel"
$ws.Range("R461").Value = "Región de O'Higgins"
$ws.Range("S461").Value = 633
$ws.Range("T461").Value = 15

# New row 462
$ws.Range("A462").Value = 11
$ws.Range("B462").Value = "Vega Monumental Concepción"
$ws.Range("C462").Value = "Bíobío"
$ws.Range("D462").Value = 45211
$ws.Range("E462").Value = 8
$ws.Range("F462").Value = "Fruta"
$ws.Range("G462").Value = 100102
$ws.Range("H462").Value = "Cítricos"
$ws.Range("I462").Value = 100102005
$ws.Range("J462").Value = "Naranja"
$ws.Range("K462").Value = "Valencia"
$ws.Range("L462").Value = "Primera"
$ws.Range("M462").Value = 200
$ws.Range("N462").Value = 9000
$ws.Range("O462").Value = 10000
$ws.Range("P462").Value = 9500
$ws.Range("Q462").Value = "$/bandeja 15 kilos granel"
$ws.Range("R462").Value = "Región de O'Higgins"
$ws.Range("S462").Value = 633
$ws.Range("T462").Value = 15
